$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the allowableKey (column B) and correctKey (column D) values for rows 2-5,
# and add the new correctKeyPress (column E) header + numeric values.

$ws.Range("B1").Value = "allowableKey"

$ws.Range("B2").Value = "['7', '4', '1','0']"
$ws.Range("D2").Value = "['7']"
$ws.Range("E1").Value = "correctKeyPress"
$ws.Range("E2").Value = 7

$ws.Range("B3").Value = "['7', '4', '1','0']"
$ws.Range("D3").Value = "['4']"
$ws.Range("E3").Value = 4

$ws.Range("B4").Value = "['7', '4', '1','0']"
$ws.Range("D4").Value = "['1']"
$ws.Range("E4").Value = 1

$ws.Range("B5").Value = "['7', '4', '1','0']"
$ws.Range("D5").Value = "['0']"
$ws.Range("E5").Value = 0

$ws.Range("E2:E5").Style = $ws.Range("D2").Style

$ws.Range("B6").Select()
